$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 767, shifting existing row 767 (and all rows
# below it) down by one, then populate the new row with the data for
# 2026/02/07.
$ws.Rows.Item(767).Insert()

$ws.Cells.Item(767, 1).NumberFormat = "@"
$ws.Cells.Item(767, 1).Value = "2026/02/07"
$ws.Cells.Item(767, 1).ClearFormats()
$ws.Cells.Item(767, 2).Value = "土"
$ws.Cells.Item(767, 3).Value = 1
$ws.Cells.Item(767, 4).Value = 75
